$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B11").Value = 9
$ws.Range("C11").Value = "RifleBullet"

$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 'Melee "Bullet"'

$ws.Range("C13").Select()
